{"js": "// Old -> new text, in document order: the worksheet date heading followed by\n// the 25 \"two digit \u00f7 one digit\" practice problems in the table.\nconst replacements = [\n  [\"2025-01-21 Tuesday\", \"2025-01-22 Wednesday\"],\n  [\"98\u00f72=\", \"88\u00f75=\"],\n  [\"42\u00f74=\", \"49\u00f79=\"],\n  [\"14\u00f73=\", \"56\u00f78=\"],\n  [\"20\u00f73=\", \"74\u00f77=\"],\n  [\"55\u00f74=\", \"51\u00f77=\"],\n  [\"49\u00f79=\", \"41\u00f73=\"],\n  [\"17\u00f72=\", \"39\u00f73=\"],\n  [\"62\u00f73=\", \"91\u00f74=\"],\n  [\"68\u00f74=\", \"54\u00f76=\"],\n  [\"43\u00f73=\", \"66\u00f79=\"],\n  [\"17\u00f75=\", \"79\u00f73=\"],\n  [\"15\u00f77=\", \"50\u00f72=\"],\n  [\"63\u00f75=\", \"25\u00f73=\"],\n  [\"23\u00f75=\", \"82\u00f74=\"],\n  [\"32\u00f77=\", \"13\u00f76=\"],\n  [\"44\u00f75=\", \"30\u00f79=\"],\n  [\"89\u00f75=\", \"72\u00f78=\"],\n  [\"73\u00f73=\", \"78\u00f76=\"],\n  [\"90\u00f75=\", \"29\u00f78=\"],\n  [\"74\u00f74=\", \"70\u00f78=\"],\n  [\"29\u00f75=\", \"88\u00f75=\"],\n  [\"81\u00f76=\", \"73\u00f75=\"],\n  [\"55\u00f72=\", \"16\u00f78=\"],\n  [\"10\u00f78=\", \"37\u00f73=\"],\n  [\"98\u00f79=\", \"83\u00f72=\"],\n];\n\n// Resolve every search range *before* mutating anything. This guarantees\n// each search matches only its original (still-unedited) text, even though\n// a couple of the new values happen to equal another entry's old value\n// (e.g. \"42\u00f74=\" -> \"49\u00f79=\" while a different cell's original text is\n// \"49\u00f79=\"), which would otherwise make a later search ambiguous.\nconst searches = replacements.map(([oldText]) =>\n  context.document.body.search(oldText, { matchCase: true })\n);\nsearches.forEach((s) => s.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < searches.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searches[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the worksheet date heading (first paragraph of the document).\n$datePara = $d.Paragraphs(1).Range\n$datePara.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark from the replace\n$datePara.Text = \"2025-01-22 Wednesday\"\n\n# Update the practice-problem table: the first table on the page, laid out as\n# groups of one data row followed by three blank spacer rows, 5 columns wide.\n# Addressing cells by (row, column) is purely positional, so it cannot be\n# confused by the fact that some of the new problem text happens to match\n# another cell's original text (e.g. \"49\u00f79=\" is both an old value and a new\n# value at different positions).\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n  1  = @(\"88\u00f75=\", \"49\u00f79=\", \"56\u00f78=\", \"74\u00f77=\", \"51\u00f77=\")\n  5  = @(\"41\u00f73=\", \"39\u00f73=\", \"91\u00f74=\", \"54\u00f76=\", \"66\u00f79=\")\n  9  = @(\"79\u00f73=\", \"50\u00f72=\", \"25\u00f73=\", \"82\u00f74=\", \"13\u00f76=\")\n  13 = @(\"30\u00f79=\", \"72\u00f78=\", \"78\u00f76=\", \"29\u00f78=\", \"70\u00f78=\")\n  17 = @(\"88\u00f75=\", \"73\u00f75=\", \"16\u00f78=\", \"37\u00f73=\", \"83\u00f72=\")\n}\n\nforeach ($row in $newValues.Keys) {\n  $rowValues = $newValues[$row]\n  for ($col = 1; $col -le $rowValues.Length; $col++) {\n    $t.Cell($row, $col).Range.Text = $rowValues[$col - 1]\n  }\n}\n"}
